# Applies the "Updated cryptos list" data refresh to Sheet1 (rows 2-51).
# Source: GitHub Actions scheduled scrape of coinranking.com, Wed Jun 28 10:08:45 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.291.18'
$ws.Range("E2").Value = '  -0.51%  '

# Row 3
$ws.Range("D3").Value = '1.859.49'
$ws.Range("E3").Value = '  -0.95%  '

# Row 4
$ws.Range("D4").Value = "'" + '0.9999'
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").Value = "'" + '233.58'

# Row 6
$ws.Range("E6").Value = '  +0.08%  '

# Row 7
$ws.Range("D7").Value = "'" + '0.4758'
$ws.Range("E7").Value = '  -1.03%  '

# Row 8
$ws.Range("D8").Value = "'" + '0.2755'
$ws.Range("E8").Value = '  -2.43%  '

# Row 9
$ws.Range("D9").Value = "'" + '0.06440'
$ws.Range("E9").Value = '  -1.32%  '

# Row 10
$ws.Range("D10").Value = '1.854.95'
$ws.Range("E10").Value = '  -1.27%  '

# Row 11
$ws.Range("D11").Value = "'" + '0.07420'
$ws.Range("E11").Value = '  -1.01%  '

# Row 12
$ws.Range("D12").Value = "'" + '16.12'
$ws.Range("E12").Value = '  -2.81%  '

# Row 13
$ws.Range("D13").Value = "'" + '4.993'
$ws.Range("E13").Value = '  -1.61%  '

# Row 14
$ws.Range("D14").Value = "'" + '85.76'
$ws.Range("E14").Value = '  -3.18%  '

# Row 15
$ws.Range("D15").Value = "'" + '0.6326'
$ws.Range("E15").Value = '  -4.58%  '

# Row 16
$ws.Range("D16").Value = '30.274.65'
$ws.Range("E16").Value = '  -0.34%  '

# Row 17
$ws.Range("D17").Value = "'" + '0.9997'
$ws.Range("E17").Value = '  -0.08%  '

# Row 18
$ws.Range("D18").Value = "'" + '12.79'
$ws.Range("E18").Value = '  -3.83%  '

# Row 19
$ws.Range("D19").Value = "'" + '230.46'
$ws.Range("E19").Value = '  +3.42%  '

# Row 20
$ws.Range("D20").Value = "'" + '0.000007369'
$ws.Range("E20").Value = '  -3.10%  '

# Row 21
$ws.Range("D21").Value = '2.097.40'
$ws.Range("E21").Value = '  -0.86%  '

# Row 22
$ws.Range("D22").Value = "'" + '1.001'
$ws.Range("E22").Value = '  +0.21%  '

# Row 23
$ws.Range("D23").Value = "'" + '5.105'
$ws.Range("E23").Value = '  -3.81%  '

# Row 24
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = "'" + '6.020'
$ws.Range("E24").Value = '  -2.77%  '

# Row 25
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = "'" + '167.97'
$ws.Range("E25").Value = '  +0.59%  '

# Row 26
$ws.Range("D26").Value = "'" + '9.271'
$ws.Range("E26").Value = '  -0.61%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = "'" + '17.91'
$ws.Range("E27").Value = '  -2.98%  '

# Row 28
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = "'" + '1.858'
$ws.Range("E28").Value = '  -5.47%  '

# Row 29
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").Value = "'" + '0.1020'
$ws.Range("E29").Value = '  +8.40%  '

# Row 30
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = "'" + '1.384'
$ws.Range("E30").Value = '  -5.23%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = "'" + '4.233'
$ws.Range("E31").Value = '  -1.61%  '

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'" + '3.910'
$ws.Range("E32").Value = '  -3.08%  '

# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = "'" + '0.04886'
$ws.Range("E33").Value = '  -2.56%  '

# Row 34
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = "'" + '1.146'
$ws.Range("E34").Value = '  -5.53%  '

# Row 35
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = "'" + '0.7232'
$ws.Range("E35").Value = '  -2.74%  '

# Row 36
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").Value = "'" + '0.9992'
$ws.Range("E36").Value = '  +0.19%  '

# Row 37
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = "'" + '2.693'
$ws.Range("E37").Value = '  -0.39%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = "'" + '0.01955'
$ws.Range("E38").Value = '  +6.98%  '

# Row 39
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = "'" + '2.633'
$ws.Range("E39").Value = '  +0.64%  '

# Row 40
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = "'" + '0.9100'
$ws.Range("E40").Value = '  +0.43%  '

# Row 41
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = "'" + '1.982'
$ws.Range("E41").Value = '  -4.00%  '

# Row 42
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = "'" + '105.77'
$ws.Range("E42").Value = '  -0.62%  '

# Row 43
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = "'" + '0.9996'
$ws.Range("E43").Value = '  -0.50%  '

# Row 44
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = "'" + '0.4119'
$ws.Range("E44").Value = '  -3.65%  '

# Row 45
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = "'" + '5.576'
$ws.Range("E45").Value = '  -4.45%  '

# Row 46
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = "'" + '7.053'
$ws.Range("E46").Value = '  -5.28%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = "'" + '61.19'
$ws.Range("E47").Value = '  -4.89%  '

# Row 48
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = "'" + '0.1208'
$ws.Range("E48").Value = '  -5.43%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'" + '8.821'
$ws.Range("E49").Value = '  -0.34%  '

# Row 50
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = "'" + '1.400'
$ws.Range("E50").Value = '  -5.03%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'" + '0.05615'
$ws.Range("E51").Value = '  -0.41%  '
